# Hour reg update for all
# Fill in the missing "week 3.3" hours for Tuesday, Wednesday and Thursday

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tuesday (row 22) - week 3.3 section
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 0

# Wednesday (row 23) - week 3.3 section
$ws.Range("E23").Value = 5
$ws.Range("G23").Value = 4.5
$ws.Range("I23").Value = 6.5

# Thursday (row 24) - week 3.3 section
$ws.Range("D24").Value = 4
$ws.Range("E24").Value = 4
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 4
$ws.Range("I24").Value = 4

# Reflect the cell that was left selected in the saved file
$ws.Range("F21").Select()

$wb.Save()
